# Apply weekly price-update shift for rows 123-190 (Berenjena sheet)
# Columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion), O (Origen),
# P (Precio $/Kg), Q (Kg o Unidades) each get the value that used to sit one
# row above, with a brand-new trailing row (190) added using what had been
# the last row's (189) values for this block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @'
123|44434|300|6000|6000|6000|$/caja 60 unidades|Región de Arica y Parinacota|100|60
124|44414|300|7000|7000|7000|$/caja 60 unidades|Región de Arica y Parinacota|117|60
125|45030|150|10000|10000|10000|$/caja 50 unidades|Región del Maule|200|50
126|44433|300|7000|7000|7000|$/caja 60 unidades|Región de Arica y Parinacota|117|60
127|44251|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
128|44299|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
129|45119|200|7000|7000|7000|$/caja 50 unidades|Región de Arica y Parinacota|140|50
130|44263|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
131|44270|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
132|44431|400|6000|6000|6000|$/caja 60 unidades|Región de Arica y Parinacota|100|60
133|44264|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
134|44991|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
135|44435|300|6000|6000|6000|$/caja 50 unidades|Región de Arica y Parinacota|120|50
136|44435|1300|6000|7000|6231|$/caja 60 unidades|Región de Arica y Parinacota|104|60
137|44588|200|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
138|45120|200|7000|7000|7000|$/caja 50 unidades|Región de Arica y Parinacota|140|50
139|44236|300|7000|7000|7000|$/caja 60 unidades|Región del Maule|117|60
140|44298|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
141|44421|200|8000|8000|8000|$/caja 60 unidades|Región de Arica y Parinacota|133|60
142|44987|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
143|44572|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
144|44245|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
145|44573|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
146|44676|150|9000|9000|9000|$/caja 50 unidades|Región del Maule|180|50
147|45002|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
148|44244|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
149|44258|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
150|44279|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
151|44424|300|8000|8000|8000|$/caja 60 unidades|Región de Arica y Parinacota|133|60
152|44249|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
153|44560|150|9000|9000|9000|$/caja 50 unidades|Región del Maule|180|50
154|44596|150|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
155|44725|300|7000|7000|7000|$/caja 50 unidades|Región de Arica y Parinacota|140|50
156|44571|200|8000|8000|8000|$/caja 60 unidades|Provincia de Chacabuco|133|60
157|44280|200|7000|7000|7000|$/caja 60 unidades|Región del Maule|117|60
158|44438|300|6000|6000|6000|$/caja 60 unidades|Región de Arica y Parinacota|100|60
159|44609|150|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
160|45062|200|8000|8000|8000|$/caja 50 unidades|Región de Arica y Parinacota|160|50
161|45008|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
162|44970|200|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
163|44242|300|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
164|44274|150|7000|7000|7000|$/caja 60 unidades|Región del Maule|117|60
165|45127|200|7000|7000|7000|$/caja 50 unidades|Región de Arica y Parinacota|140|50
166|45131|200|8000|8000|8000|$/caja 50 unidades|Región de Arica y Parinacota|160|50
167|44622|180|6000|6000|6000|$/caja 50 unidades|Región del Maule|120|50
168|44273|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
169|44554|100|10000|10000|10000|$/caja 50 unidades|Región del Maule|200|50
170|44272|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
171|44589|200|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
172|44565|150|9000|9000|9000|$/caja 50 unidades|Región del Maule|180|50
173|45132|200|7000|7000|7000|$/caja 50 unidades|Región de Arica y Parinacota|140|50
174|45111|300|6000|6000|6000|$/caja 50 unidades|Región de Arica y Parinacota|120|50
175|44417|300|7000|7000|7000|$/caja 60 unidades|Región de Arica y Parinacota|117|60
176|44664|200|9000|9000|9000|$/caja 50 unidades|Región del Maule|180|50
177|44468|200|8000|8000|8000|$/caja 50 unidades|Región de Arica y Parinacota|160|50
178|44253|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
179|44235|200|7000|7000|7000|$/caja 60 unidades|Región del Maule|117|60
180|44239|200|8000|8000|8000|$/caja 60 unidades|Región del Maule|133|60
181|44971|200|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
182|44985|200|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
183|44729|300|8000|8000|8000|$/caja 50 unidades|Región de Arica y Parinacota|160|50
184|45112|300|6000|6000|6000|$/caja 50 unidades|Región de Arica y Parinacota|120|50
185|45112|300|6000|6000|6000|$/caja 50 unidades|Región de Arica y Parinacota|120|50
186|44574|150|8000|8000|8000|$/caja 50 unidades|Región del Maule|160|50
187|44566|200|9000|9000|9000|$/caja 50 unidades|Región del Maule|180|50
188|44594|200|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
189|45121|200|7000|7000|7000|$/caja 50 unidades|Región de Arica y Parinacota|140|50
190|44607|300|7000|7000|7000|$/caja 50 unidades|Región del Maule|140|50
'@

$lines = $data -split "`n"

foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $dVal = [double]$parts[1]
    $jVal = [double]$parts[2]
    $kVal = [double]$parts[3]
    $lVal = [double]$parts[4]
    $mVal = [double]$parts[5]
    $nVal = $parts[6]
    $oVal = $parts[7]
    $pVal = [double]$parts[8]
    $qVal = [double]$parts[9]

    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($row, 10).Value = $jVal
    $ws.Cells.Item($row, 11).Value = $kVal
    $ws.Cells.Item($row, 12).Value = $lVal
    $ws.Cells.Item($row, 13).Value = $mVal
    $ws.Cells.Item($row, 14).Value = $nVal
    $ws.Cells.Item($row, 15).Value = $oVal
    $ws.Cells.Item($row, 16).Value = $pVal
    $ws.Cells.Item($row, 17).Value = $qVal
}

# Row 190 is brand new - fill in the columns that stay constant across the
# whole Berenjena block (same as every other data row).
$ws.Cells.Item(190, 1).Value = 5
$ws.Cells.Item(190, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(190, 3).Value = "Maule"
$ws.Cells.Item(190, 5).Value = 7
$ws.Cells.Item(190, 6).Value = 100112001
$ws.Cells.Item(190, 7).Value = "Berenjena"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 18).Value = "Hortaliza"
